# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    carrying the same per-fund breakdown layout as the other quarter sheets.
# 2. Insert a new leading row into "总计" summarising the new quarter, pushing
#    the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$wsQ4    = $wb.Worksheets.Item("2021-Q4")
$wsTotal = $wb.Worksheets.Item("总计")

# --- 1. New "2022-Q1" sheet, inserted just before "总计" ------------------

# Duplicate the 2021-Q4 sheet (placed right before "总计") so the new sheet
# naturally inherits its layout: sheetPr/outline settings, page margins and
# the bold/centered/thin-bordered header + index-column formatting, instead
# of re-building all of that cell by cell.
$wsQ4.Copy($wsTotal)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# 2021-Q4 has 3 data rows; 2022-Q1 only has 1 - drop the extra two.
$newSheet.Rows.Item(3).Resize(2).Delete()

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'003981"
$newSheet.Range("C2").Value = "中银证券瑞益灵活配置混合C"
$newSheet.Range("D2").Value = "'0.21"
$newSheet.Range("E2").Value = "'89.21"
$newSheet.Range("F2").Value = "'3.09"
$newSheet.Range("G2").Value = "'0.0065"
$newSheet.Range("H2").Value = 9

# --- 2. New leading row in "总计" for the 2022-Q1 summary ------------------

# Re-resolve "总计" by name: inserting the sheet above shifted its position,
# and the earlier $wsTotal handle tracks that *position*, not the sheet.
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("A2:D2").Insert(-4121)

# The inserted row inherits the header row's (bold/bordered) formatting;
# reset the plain data cells and re-apply the index-column style (copied
# from the row below, which still carries it) so the new row matches its
# siblings (B:D unstyled, A styled).
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.01

# The inserted row pushed the old index column (0,1,2) down one row without
# renumbering it; re-sequence it (0..3) to keep it a contiguous row index.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
